# Horarios Línea 141 - update scrape snapshot (05:47:29 -> 06:15:23) with
# refreshed row counts/minute offsets and newly scraped rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "LP1912": dimension A1:E21 -> A1:E27
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 06:15:23"
$ws1.Range("A3").Value = "Total filas: 22"

$rows1 = @(
    @(8,  "06:15:23", "06:16", "215A_EL PATO",  1,   "LP1912"),
    @(9,  "06:15:23", "06:30", "23_HERNANDEZ",  15,  "LP1912"),
    @(10, "06:15:23", "06:34", "11_ETCHEVERRY", 19,  "LP1912"),
    @(11, "06:15:23", "06:39", "17X38_ROMERO",  24,  "LP1912"),
    @(12, "06:15:23", "06:41", "16_SANTA ANA",  26,  "LP1912"),
    @(13, "06:15:23", "06:57", "215A_EL PATO",  42,  "LP1912"),
    @(14, "06:15:23", "06:59", "225_GOMEZ",     44,  "LP1912"),
    @(16, "06:15:23", "07:16", "215C_EL PATO",  61,  "LP1912"),
    @(17, "06:15:23", "07:19", "14_ABASTO",     64,  "LP1912"),
    @(18, "06:15:23", "07:21", "16_SANTA ANA",  66,  "LP1912"),
    @(19, "06:15:23", "07:21", "23_HERNANDEZ",  66,  "LP1912"),
    @(20, "06:15:23", "07:29", "17X38_ROMERO",  74,  "LP1912"),
    @(21, "06:15:23", "07:35", "10_OLMOS",      80,  "LP1912"),
    @(22, "05:47:29", "07:36", "27_EL RETIRO",  109, "LP1912"),
    @(23, "06:15:23", "07:37", "27_EL RETIRO",  82,  "LP1912"),
    @(24, "06:15:23", "07:55", "14_ABASTO",     100, "LP1912"),
    @(25, "06:15:23", "08:00", "17_ROMERO",     105, "LP1912"),
    @(26, "06:15:23", "08:11", "10_OLMOS",      116, "LP1912"),
    @(27, "06:15:23", "08:13", "15X38_ABASTO",  118, "LP1912")
)

foreach ($r in $rows1) {
    $row = $r[0]
    $ws1.Cells.Item($row, 1).Value = $r[1]
    $ws1.Cells.Item($row, 2).Value = $r[2]
    $ws1.Cells.Item($row, 3).Value = $r[3]
    $ws1.Cells.Item($row, 4).Value = $r[4]
    $ws1.Cells.Item($row, 5).Value = $r[5]
}

# ---------------------------------------------------------------------------
# Sheet "LP1912-215": dimension A1:E8 -> A1:E9
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 06:15:23"
$ws2.Range("A3").Value = "Total filas: 4"

$rows2 = @(
    @(6, "06:15:23", "06:16", "215A_EL PATO", 1,  "LP1912"),
    @(7, "06:15:23", "06:57", "215A_EL PATO", 42, "LP1912"),
    @(9, "06:15:23", "07:16", "215C_EL PATO", 61, "LP1912")
)

foreach ($r in $rows2) {
    $row = $r[0]
    $ws2.Cells.Item($row, 1).Value = $r[1]
    $ws2.Cells.Item($row, 2).Value = $r[2]
    $ws2.Cells.Item($row, 3).Value = $r[3]
    $ws2.Cells.Item($row, 4).Value = $r[4]
    $ws2.Cells.Item($row, 5).Value = $r[5]
}

# ---------------------------------------------------------------------------
# Sheet "6203-6173": dimension unchanged (A1:E6)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 06:15:23"
$ws3.Cells.Item(6, 1).Value = "06:15:23"
$ws3.Cells.Item(6, 4).Value = 88
